# Apply Price (D) / Volume(1h) (E) updates from the latest cryptos refresh.
# Numeric-looking Price strings must stay as TEXT (matching the original
# inlineStr cells), so they get an apostrophe prefix plus a Style reset
# back to "Normal" to avoid leaving a stray quote-prefixed number format.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.319.01"
$ws.Range("E2").Value = "  -3.14%  "

$ws.Range("D3").Value = "1.936.59"
$ws.Range("E3").Value = "  -3.25%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").Value = "'250.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.00%  "

$ws.Range("E6").Value = "  -7.00%  "

$ws.Range("D7").Value = "'1.001"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.15%  "

$ws.Range("D8").Value = "'0.3315"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.77%  "

$ws.Range("D9").Value = "'27.77"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.50%  "

$ws.Range("D10").Value = "'0.07291"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.76%  "

$ws.Range("D11").Value = "'0.8106"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.65%  "

$ws.Range("D12").Value = "'0.08097"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.30%  "

$ws.Range("D13").Value = "1.937.16"
$ws.Range("E13").Value = "  -3.27%  "

$ws.Range("D14").Value = "'5.504"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.74%  "

$ws.Range("D15").Value = "'94.69"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -6.44%  "

$ws.Range("D16").Value = "'15.16"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.26%  "

$ws.Range("D17").Value = "30.334.21"
$ws.Range("E17").Value = "  -3.04%  "

$ws.Range("D18").Value = "'0.000008323"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.91%  "

$ws.Range("D19").Value = "'251.45"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -7.95%  "

$ws.Range("D20").Value = "'5.899"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.04%  "

$ws.Range("D21").Value = "2.191.61"
$ws.Range("E21").Value = "  -2.85%  "

$ws.Range("D22").Value = "'1.001"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.18%  "

$ws.Range("E23").Value = "  +0.07%  "

$ws.Range("D24").Value = "'7.007"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.97%  "

$ws.Range("D25").Value = "'9.767"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.61%  "

$ws.Range("D26").Value = "'163.51"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.57%  "

$ws.Range("D27").Value = "'2.390"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.55%  "

$ws.Range("D28").Value = "'19.28"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.69%  "

$ws.Range("D29").Value = "'0.1319"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -7.15%  "

$ws.Range("D30").Value = "'1.566"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.18%  "

$ws.Range("E31").Value = "  -1.48%  "

$ws.Range("D32").Value = "'4.431"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.00%  "

$ws.Range("D33").Value = "'4.187"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.87%  "

$ws.Range("D34").Value = "'0.05204"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.18%  "

$ws.Range("D35").Value = "'1.291"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.93%  "

$ws.Range("D36").Value = "'0.7504"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.24%  "

$ws.Range("D37").Value = "'2.743"
$ws.Range("D37").Style = "Normal"

$ws.Range("D38").Value = "'0.01978"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.81%  "

$ws.Range("D39").Value = "'2.822"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.72%  "

$ws.Range("D40").Value = "'79.40"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -8.80%  "

$ws.Range("D41").Value = "'6.357"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.83%  "

$ws.Range("D42").Value = "'0.4539"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.50%  "

$ws.Range("D43").Value = "'2.026"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.37%  "

$ws.Range("D44").Value = "'0.8454"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.58%  "

$ws.Range("D45").Value = "'1.000"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.14%  "

$ws.Range("D46").Value = "'101.74"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.75%  "

$ws.Range("D47").Value = "'9.780"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.06%  "

$ws.Range("D48").Value = "'7.474"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.34%  "

$ws.Range("D49").Value = "'36.82"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.88%  "

$ws.Range("D50").Value = "'0.4194"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.99%  "

$ws.Range("D51").Value = "'0.06035"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.36%  "

